# Generate Report for Handback
# Refresh the generated timestamps on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for eb8bbaed-96ec-...md
$wsOverview.Range("G2").Value = "2016-08-26 19:08:47"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-26 19:08:42"
$wsZhCn.Range("K2").Value = "2016-08-26 19:08:59"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-26 19:08:47"
$wsDeDe.Range("K2").Value = "2016-08-26 19:09:12"
